# Sprint Burndown Chart Template — apply commit edits
#
# Summary of the change (per the OOXML diff):
#   1. E18 (row 18, "Day 9" column) gets a new logged value of 1 hour.
#      This ripples through the shared "Result" formula in column P
#      (P18 flips from "-" to "Incomplete") and through the Remaining /
#      Ideal-Trend rollup rows 26-27, which in turn re-seed the chart's
#      cached numCache points — all of that is handled automatically by
#      the engine's recalculation, we only need to author the literal
#      edit cell.
#   2. The data-validation rule that used to cover E6:O27 is narrowed to
#      F6:O27 + E26:E27, and a brand new rule (0-24, instead of -24-24)
#      is added for E6:E25 — exactly what Excel does natively when you
#      apply a new validation on top of part of an existing validated
#      range.
#   3. The sheet's view scrolled/selected a different cell: selection
#      moves from K22 to F20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Log 1 hour of work against "Day 9" (row 18) -----------------------
$ws.Range("E18").Value = 1

# --- 2. Data validation: applying a brand new (0..24) validation rule on
#        top of the previously-validated E6:E25 sub-range makes Excel
#        itself shrink the old (-24..24) rule's sqref down to the cells
#        that are still only covered by it (F6:O27 plus the E26:E27
#        tail), while the new rule takes E6:E25. ---------------------------
$ws.Range("E6:E25").Validation.Delete()
$ws.Range("E6:E25").Validation.Add(2, 1, 1, 0, 24)

# --- 3. Move the active selection to F20 -----------------------------------
$ws.Range("F20").Select() | Out-Null
